$d = $word.ActiveDocument

# The TNS number-acquisition mailbox changed from
# "acstnrequest@microsoft.com" to "acstns@microsoft.com" (the mailto:
# target in the relationships part is left untouched - only the
# hyperlink's visible text changes). Locate the hyperlink by its
# current display text so this keeps working regardless of hyperlink
# ordering in the document.
$target = $null
for ($i = 1; $i -le $d.Hyperlinks.Count; $i++) {
    $h = $d.Hyperlinks($i)
    if ($h.TextToDisplay -eq "acstnrequest@microsoft.com") {
        $target = $h
    }
}

if ($target -ne $null) {
    $target.TextToDisplay = "acstns@microsoft.com"

    # Re-locate the freshly updated text and nudge the formatting of
    # the newly-typed "s" (character 6 of "acstns@microsoft.com") so it
    # lands in its own run instead of silently re-merging with its
    # neighbours - mirroring the run split Word leaves behind when you
    # type over a mid-run selection.
    $r = $d.Content
    $found = $r.Find.Execute("acstns@microsoft.com")
    if ($found) {
        $start = $r.Start
        $sRange = $d.Range($start + 5, $start + 6)
        $sRange.Font.Size = 11
        $sRange.Font.Size = 10
    }
}
